$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-11 Wednesday" "2024-12-12 Thursday"

Replace-Text "44×28=" "74×40="
Replace-Text "78×28=" "23×24="
Replace-Text "22×95=" "87×18="
Replace-Text "41×34=" "73×75="
Replace-Text "19×23=" "83×78="
Replace-Text "66×23=" "94×17="
Replace-Text "39×26=" "40×79="
Replace-Text "33×75=" "13×91="
Replace-Text "98×56=" "86×60="
Replace-Text "79×21=" "39×28="
Replace-Text "48×69=" "34×86="
Replace-Text "99×37=" "85×95="
Replace-Text "81×53=" "13×88="
Replace-Text "88×95=" "79×50="
Replace-Text "55×82=" "97×30="
Replace-Text "18×25=" "40×67="
Replace-Text "51×77=" "32×95="
Replace-Text "13×32=" "59×66="
Replace-Text "48×64=" "50×11="
Replace-Text "78×11=" "97×94="
Replace-Text "64×46=" "81×69="
Replace-Text "90×54=" "52×43="
Replace-Text "75×62=" "16×76="
Replace-Text "57×74=" "21×58="
Replace-Text "73×63=" "94×59="
